$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "big and tall running pants for men",
    "girls compression knee sleeve",
    "football pants with pads",
    "compression padded shorts",
    "spandex for men pants",
    "elastic knee pad",
    "knee pads gel construction",
    "capri shorts for men",
    "football padded pants",
    "girls volleyball kneepads",
    "black capris men",
    "womens basketball pants",
    "baseball material",
    "youth sports leggings",
    "protect knee pads",
    "knees pad",
    "black youth knee pads",
    "exercise pads for knees",
    "gym shorts above knee for men",
    "knee sleeve basketball youth",
    "knee pads construction",
    "work knee pad",
    "bjj knee",
    "knee sleeve bjj",
    "knee pads under",
    "mens long cycling pants",
    "mens pad",
    "boys long baseball pants",
    "mens basketball gear",
    "girl sliding shorts",
    "calf silicone pads",
    "compression knee sleeve men pair",
    "girls black softball pants",
    "hip protector pads",
    "flexible work pants for men",
    "knee protection pads",
    "lightweight work pants for men",
    "youth baseball pants long",
    "knees pads work",
    "tights boys",
    "mens shorts long below knee",
    "knee sleeve youth",
    "snowboarding pants men",
    "baseball shorts for men",
    "baseball mens pants",
    "knee compression sleeve - reduce strain & swelling",
    "pads men",
    "basketball sleeve youth leg",
    "thigh pads football",
    "compression volleyball",
    "leggings for mens",
    "mens yoga pants",
    "padded football pants",
    "spandex capris",
    "water knee hockey",
    "compression pants sleeves",
    "knee sleeve padded",
    "knees pads for construction",
    "tight capri",
    "mens baseball compression shorts",
    "mens running knee compression",
    "black football leggings",
    "knee sleeves basketball youth",
    "cycling knee pads",
    "construction knee pad",
    "compression calf leggings",
    "baseball youth compression sleeve",
    "6 pairs of leggings",
    "basketball shorts for men pack of 5",
    "compression pants youth boys",
    "impact shorts men",
    "large knee pad",
    "baseball compression sleeve",
    "boys sports leggings",
    "volleyball spandex pack",
    "baseball pants youth large",
    "boys paintball pants",
    "yoga pants mens",
    "calf tear compression sleeve",
    "compression shorts men long length",
    "yoga hand pads",
    "knee sleeves with padding",
    "athletic capri leggings",
    "pants compression men",
    "basketball padding",
    "knee pads for men floor work",
    "youth knee sleeve wrestling",
    "professional construction knee pads",
    "youth basketball",
    "basketball compression knee sleeve",
    "black softball pants youth girls",
    "hex gear wash",
    "knee construction pads",
    "youth girls softball pants",
    "compression tight pants",
    "male workout leggings",
    "boys xl baseball pants",
    "thick leggings for men",
    "knee pads for",
    "adult football girdle"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $newValues[$i]
}
